# Update the "Unemployment Rate" column (P) on the "Quarterly Data Levels"
# sheet so each quarter's value is rounded to one decimal place instead of
# the raw repeating-decimal average, and backfill the newly-available
# "Employment level" figure (N) for the latest quarter (row 124).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> new (rounded) Unemployment Rate value for column P
$pUpdates = @{
    14  = 8.8
    16  = 7.3
    17  = 6.7
    18  = 6.3
    20  = 5.4
    24  = 4.2
    26  = 4
    28  = 4.4
    29  = 4.4
    30  = 4.5
    31  = 4.6
    33  = 5.3
    34  = 4.8
    35  = 4.8
    37  = 4.8
    38  = 5.1
    39  = 4.7
    42  = 4.4
    43  = 4.9
    45  = 4.5
    47  = 4.8
    49  = 4.6
    54  = 5.4
    55  = 5.9
    56  = 7.4
    57  = 8.6
    58  = 10.9
    61  = 13.8
    63  = 14.4
    65  = 15.7
    66  = 15.1
    67  = 15.1
    69  = 15.9
    70  = 15.9
    71  = 15.5
    74  = 14.5
    75  = 14.4
    76  = 13.3
    77  = 12.9
    78  = 12.8
    80  = 11.7
    81  = 10.9
    82  = 10.6
    83  = 10.1
    84  = 9.6
    85  = 9.4
    86  = 9
    88  = 8.2
    89  = 7.5
    90  = 7.3
    91  = 6.6
    93  = 6.4
    94  = 5.9
    95  = 5.9
    97  = 5.7
    98  = 5
    99  = 5.2
    101 = 4.8
    102 = 4.9
    103 = 5.1
    105 = 6.3
    106 = 7.3
    107 = 7
    108 = 5.5
    109 = 5.2
    110 = 4.9
    112 = 4.2
    113 = 4.4
    115 = 4.2
    117 = 4.5
    118 = 4.3
    120 = 4.2
    122 = 4.5
    124 = 5
}

foreach ($row in $pUpdates.Keys) {
    $ws.Range("P$row").Value = $pUpdates[$row]
}

# Newly available Employment level figure for the latest quarter.
$ws.Range("N124").Value = 2825.5
